$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "John"
$ws.Range("B4").Value = "S5678332G"
$ws.Range("C4").Value = 28
$ws.Range("D4").Value = "Married"
$ws.Range("E4").Value = "password"
[void]$ws.Range("E5").Select()
